$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.876.63"
Set-TextValue $ws.Range("E2") "  +0.40%  "
Set-TextValue $ws.Range("D3") "3.556.14"
Set-TextValue $ws.Range("E3") "  -1.14%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  -0.08%  "
Set-TextValue $ws.Range("D5") "609.78"
Set-TextValue $ws.Range("E5") "  +0.05%  "
Set-TextValue $ws.Range("D6") "145.65"
Set-TextValue $ws.Range("E6") "  -1.99%  "
Set-TextValue $ws.Range("D7") "3.555.60"
Set-TextValue $ws.Range("E7") "  -1.09%  "
Set-TextValue $ws.Range("E8") "  -0.12%  "
Set-TextValue $ws.Range("E9") "  +5.79%  "
Set-TextValue $ws.Range("E10") "  -1.89%  "
Set-TextValue $ws.Range("D11") "7.82"
Set-TextValue $ws.Range("E11") "  -2.83%  "
Set-TextValue $ws.Range("E12") "  +0.79%  "
Set-TextValue $ws.Range("D13") "4.160.90"
Set-TextValue $ws.Range("E13") "  -1.16%  "
Set-TextValue $ws.Range("D14") "0.0000199"
Set-TextValue $ws.Range("E14") "  -5.17%  "
Set-TextValue $ws.Range("D15") "29.30"
Set-TextValue $ws.Range("E15") "  -1.78%  "
Set-TextValue $ws.Range("D16") "3.553.55"
Set-TextValue $ws.Range("E16") "  -0.92%  "
Set-TextValue $ws.Range("B17") "WrappedBTC"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "66.708.81"
Set-TextValue $ws.Range("E17") "  +0.00%  "
Set-TextValue $ws.Range("B18") "TRON"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D18") "0.117"
Set-TextValue $ws.Range("E18") "  +0.88%  "
Set-TextValue $ws.Range("E19") "  -3.80%  "
Set-TextValue $ws.Range("D20") "6.29"
Set-TextValue $ws.Range("E20") "  -1.18%  "
Set-TextValue $ws.Range("E21") "  -2.03%  "
Set-TextValue $ws.Range("D22") "428.99"
Set-TextValue $ws.Range("E22") "  +0.35%  "
Set-TextValue $ws.Range("D23") "0.603"
Set-TextValue $ws.Range("E23") "  -2.31%  "
Set-TextValue $ws.Range("D24") "77.86"
Set-TextValue $ws.Range("E24") "  -1.26%  "
Set-TextValue $ws.Range("D25") "3.698.87"
Set-TextValue $ws.Range("E25") "  -1.18%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  -0.06%  "
Set-TextValue $ws.Range("E27") "  -5.14%  "
Set-TextValue $ws.Range("E28") "  -3.29%  "
Set-TextValue $ws.Range("E29") "  -1.26%  "
Set-TextValue $ws.Range("D30") "9.13"
Set-TextValue $ws.Range("E30") "  -2.67%  "
Set-TextValue $ws.Range("E31") "  +0.05%  "
Set-TextValue $ws.Range("E32") "  -0.09%  "
Set-TextValue $ws.Range("D33") "3.565.49"
Set-TextValue $ws.Range("E33") "  -0.82%  "
Set-TextValue $ws.Range("E34") "  -3.43%  "
Set-TextValue $ws.Range("E35") "  +0.00%  "
Set-TextValue $ws.Range("E36") "  -7.81%  "
Set-TextValue $ws.Range("E37") "  -1.79%  "
Set-TextValue $ws.Range("E38") "  -3.06%  "
Set-TextValue $ws.Range("D39") "177.86"
Set-TextValue $ws.Range("E39") "  +0.05%  "
Set-TextValue $ws.Range("D40") "5.34"
Set-TextValue $ws.Range("E40") "  -5.52%  "
Set-TextValue $ws.Range("D41") "0.0835"
Set-TextValue $ws.Range("E41") "  -2.55%  "
Set-TextValue $ws.Range("D42") "5.07"
Set-TextValue $ws.Range("E42") "  -3.42%  "
Set-TextValue $ws.Range("E43") "  -3.63%  "
Set-TextValue $ws.Range("D44") "45.60"
Set-TextValue $ws.Range("E44") "  -1.34%  "
Set-TextValue $ws.Range("E45") "  -5.76%  "
Set-TextValue $ws.Range("D46") "0.999"
Set-TextValue $ws.Range("E46") "  -0.06%  "
Set-TextValue $ws.Range("D47") "2.42"
Set-TextValue $ws.Range("E47") "  -5.90%  "
Set-TextValue $ws.Range("E48") "  -0.12%  "
Set-TextValue $ws.Range("D49") "23.51"
Set-TextValue $ws.Range("E49") "  -2.76%  "
Set-TextValue $ws.Range("E50") "  -4.56%  "
Set-TextValue $ws.Range("E51") "  -2.85%  "
